$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: refresh the existing weather-station observation window ---
$ws.Range("B2").Value = 45638.32013888889
$ws.Range("C2").Value = 45639.31736111111
$ws.Range("D2").Value = -0.7
$ws.Range("E2").Value = 5.5
$ws.Range("F2").Value = 1.8
$ws.Range("G2").Value = 1.3

# --- Row 3: new weather-station observation window ---
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "79049004"
$ws.Range("B3").Value = 45639.31736111111
$ws.Range("C3").Value = 45640.4375
$ws.Range("D3").Value = -0.6
$ws.Range("E3").Value = 2.8
$ws.Range("F3").Value = 1.72
$ws.Range("G3").Value = 2

# --- Row 4: new weather-station observation window ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "79049004"
$ws.Range("B4").Value = 45640.4375
$ws.Range("C4").Value = 45641.34861111111
$ws.Range("D4").Value = 1.5
$ws.Range("E4").Value = 7.1
$ws.Range("F4").Value = 4.51
$ws.Range("G4").Value = 4.6

# Match the date / numeric display formats already used on row 2
$ws.Range("B3:C4").NumberFormat = $ws.Range("B2").NumberFormat
$ws.Range("D3:G4").NumberFormat = $ws.Range("D2").NumberFormat
